# Apply the text updates described by the diff (v1.0.3 -> v1.1).
# Two shared strings change, and they are reused across 4 test cases in the
# sheet (rows 11/12, 23/24, 34/35, 53/54), each following the same pattern:
#   column D of the "Indica uma prestação..." row  -> updated SYSTEM text
#   column B of the next row ("Clica para analisar...") -> updated Chefe text

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$oldExpected = "SYSTEM Exibe os detalhes relativos àquela prestação de contas (dados básicos da solicitação e documentos anexos); Exibe o histórico da tramitação da prestação de contas."
$newExpected = "SYSTEM Exibe os detalhes relativos àquela prestação de contas (nome do beneficiário, dados básicos da solicitação e documentos anexos); Exibe o histórico da tramitação da prestação de contas."

$oldStep = "Chefe Clica para analisar a prestação de contas."
$newStep = "Chefe Verifica o histório da tramitação da prestação de contas e clica para analisar a prestação de contas."

$dCells = @("D11", "D23", "D34", "D53")
foreach ($cellRef in $dCells) {
    $cell = $ws.Range($cellRef)
    if ($cell.Value2 -eq $oldExpected) {
        $cell.Value2 = $newExpected
    }
}

$bCells = @("B12", "B24", "B35", "B54")
foreach ($cellRef in $bCells) {
    $cell = $ws.Range($cellRef)
    if ($cell.Value2 -eq $oldStep) {
        $cell.Value2 = $newStep
    }
}
